# Refresh cryptocurrency price/volume snapshot (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '40.112.59'
$ws.Range('E2').Value = '  -1.83%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.345.21'
$ws.Range('E3').Value = '  -3.11%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.20%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D5').Value = '310.47'
$ws.Range('E5').Value = '  -1.66%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D6').Value = '85.48'
$ws.Range('E6').Value = '  -3.69%  '

# Row 7: XRP
$ws.Range('E7').Value = '  -1.48%  '

# Row 8: USDC
$ws.Range('E8').Value = '  -0.06%  '

# Row 9: Cardano
$ws.Range('D9').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D9').Value = '0.483'
$ws.Range('E9').Value = '  -2.24%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  -2.06%  '

# Row 11: Avalanche
$ws.Range('D11').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D11').Value = '30.02'
$ws.Range('E11').Value = '  -5.70%  '

# Row 12: TRON
$ws.Range('E12').Value = '  +1.10%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range('D13').Value = '2.707.47'
$ws.Range('E13').Value = '  -3.20%  '

# Row 14: Polkadot
$ws.Range('D14').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D14').Value = '6.42'
$ws.Range('E14').Value = '  -3.83%  '

# Row 15: Chainlink
$ws.Range('D15').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D15').Value = '14.73'
$ws.Range('E15').Value = '  -5.81%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '2.365.09'
$ws.Range('E16').Value = '  -2.39%  '

# Row 17: Polygon
$ws.Range('E17').Value = '  -0.85%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '40.089.07'
$ws.Range('E18').Value = '  -2.02%  '

# Row 19: ShibaInu
$ws.Range('E19').Value = '  -1.65%  '

# Row 20: Uniswap
$ws.Range('E20').Value = '  -1.91%  '

# Row 21: Litecoin
$ws.Range('D21').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D21').Value = '68.02'
$ws.Range('E21').Value = '  -4.85%  '

# Row 22: InternetComputer(DFINITY)
$ws.Range('D22').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D22').Value = '10.70'
$ws.Range('E22').Value = '  -2.42%  '

# Row 23: BitcoinCash
$ws.Range('D23').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D23').Value = '235.36'
$ws.Range('E23').Value = '  +0.21%  '

# Row 24: PancakeSwap
$ws.Range('E24').Value = '  -4.70%  '

# Row 25: Dai
$ws.Range('E25').Value = '  -0.09%  '

# Row 26: ImmutableX
$ws.Range('D26').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D26').Value = '1.82'
$ws.Range('E26').Value = '  -2.22%  '

# Row 27: EthereumClassic
$ws.Range('D27').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D27').Value = '23.60'
$ws.Range('E27').Value = '  -1.50%  '

# Row 28: Toncoin
$ws.Range('E28').Value = '  -3.51%  '

# Row 29: Cosmos
$ws.Range('E29').Value = '  -2.59%  '

# Row 30: InjectiveProtocol
$ws.Range('E30').Value = '  -0.32%  '

# Row 31: Monero
$ws.Range('D31').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D31').Value = '153.46'
$ws.Range('E31').Value = '  -1.41%  '

# Row 32: FirstDigitalUSD
$ws.Range('E32').Value = '  -0.21%  '

# Row 33: Filecoin
$ws.Range('E33').Value = '  -2.33%  '

# Row 34: WEMIXToken
$ws.Range('E34').Value = '  -2.65%  '

# Row 35: Hedera
$ws.Range('E35').Value = '  -2.91%  '

# Row 36: Stellar
$ws.Range('E36').Value = '  -0.57%  '

# Row 37: LidoDAOToken
$ws.Range('E37').Value = '  -3.96%  '

# Row 38: Kaspa
$ws.Range('D38').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D38').Value = '0.0982'
$ws.Range('E38').Value = '  -1.52%  '

# Row 39: Celestia
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D39').Value = '1.72'
$ws.Range('E39').Value = '  -2.52%  '

# Row 40: ARBITRUM
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D40').Value = '15.53'
$ws.Range('E40').Value = '  -6.21%  '

# Row 41: RenderToken
$ws.Range('E41').Value = '  +0.39%  '

# Row 42: Maker
$ws.Range('D42').Value = '1.966.26'
$ws.Range('E42').Value = '  -1.04%  '

# Row 43: ApeXProtocol
$ws.Range('D43').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D43').Value = '2.24'
$ws.Range('E43').Value = '  -2.03%  '

# Row 44: EnergySwap
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D44').Value = '0.0265'
$ws.Range('E44').Value = '  -3.13%  '

# Row 45: VeChain
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D45').Value = '17.62'
$ws.Range('E45').Value = '  -5.26%  '

# Row 46: FraxShare
$ws.Range('D46').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D46').Value = '9.33'
$ws.Range('E46').Value = '  -0.68%  '

# Row 47: NEARProtocol
$ws.Range('E47').Value = '  -5.65%  '

# Row 48: RocketPoolETH
$ws.Range('D48').Value = '2.569.89'
$ws.Range('E48').Value = '  -3.28%  '

# Row 49: Aave
$ws.Range('D49').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D49').Value = '93.09'
$ws.Range('E49').Value = '  -1.26%  '

# Row 50: BitcoinSV
$ws.Range('D50').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D50').Value = '70.49'
$ws.Range('E50').Value = '  -3.02%  '

# Row 51: MultiversX
$ws.Range('D51').NumberFormat = '@'  # keep as text, not a number
$ws.Range('D51').Value = '50.37'
$ws.Range('E51').Value = '  -2.38%  '
